$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15: "VGA Control" task entry, copied from the last existing row (14)
# so it inherits the same number formats / alignment, then the actual
# values for the new entry are written on top.
$ws.Range("A14:F14").Copy($ws.Range("A15:F15"))

$ws.Range("A15").Value = "30.3.2020"
$ws.Range("B15").Value = 0.52083333333333337
$ws.Range("C15").Value = 0.53125
$ws.Range("D15").Formula = "=C15-B15"
$ws.Range("E15").Value = "VGA Control"
$ws.Range("F15").Value = "Create necessary files"

# Reflect where the cursor ends up after typing the new row (matches the
# saved selection in the workbook).
[void]$ws.Range("A16").Select()
